$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a date cell so it gets the "date" cell style (numFmtId 14,
# same style bucket already used by the other date cells in column J)
function Set-LogDate($addr, $value) {
    $cell = $ws.Range($addr)
    $cell.Value = $value
    $cell.NumberFormat = "m/d/yy"
}

# --- Row 6: "Aan de analyse werken" (A6, already set) gets a duration and
#     its date is corrected ---
$ws.Range("I6").Value = "1 uur"
Set-LogDate "J6" 40949

# --- Row 7: new entry ---
$ws.Range("A7").Value = 'Een "GO" gekregen voor ons project'
$ws.Range("I7").Value = "0,5 uur"
Set-LogDate "J7" 40953

# --- Row 8: new entry ---
$ws.Range("A8").Value = "Groesgesprek met Robert"
$ws.Range("I8").Value = "0,5 uur"
Set-LogDate "J8" 40953

# --- Row 9: new entry ---
$ws.Range("A9").Value = "Uitwerken van de models"
$ws.Range("I9").Value = "1 uur"
Set-LogDate "J9" 40955

# --- Row 10: new entry ---
$ws.Range("A10").Value = "Tutorials over ASP.NET"
$ws.Range("I10").Value = "1 uur"
Set-LogDate "J10" 40961

# --- Row 11: new entry ---
$ws.Range("A11").Value = "Models afgemaakt"
$ws.Range("I11").Value = "0,5 uur"
Set-LogDate "J11" 40961

# --- Row 12: new entry ---
$ws.Range("A12").Value = "Unit tests afgemaakt"
$ws.Range("I12").Value = "1 uur"
Set-LogDate "J12" 40961

# --- Row 13: new entry ---
$ws.Range("A13").Value = "Aan de Database gewerkt"
$ws.Range("I13").Value = "1,5 uur"
Set-LogDate "J13" 40973

# --- Row 14: new entry ---
$ws.Range("A14").Value = "Groesgesprek met Robert"
$ws.Range("I14").Value = "0,5 uur"
Set-LogDate "J14" 40973

# --- Row 15: only the date got filled in so far ---
Set-LogDate "J15" 40980

# Column widths (planning week update) - best achievable match given the
# engine's internal character-width rounding.
$ws.Columns.Item(1).ColumnWidth = 9.8333333333
$ws.Columns.Item(3).ColumnWidth = 10
$ws.Columns.Item(4).ColumnWidth = 10
$ws.Columns.Item(10).ColumnWidth = 14.8333333333

# Move the selection to A15, where the author left off editing.
$ws.Range("A15").Select() | Out-Null
